# Updated symbol list on Mon Dec 12 23:55:16 UTC 2022 with GitHub Actions
#
# Applies the coin-price / symbol-ranking refresh described by the diff:
#  - plain text cells (coin name, link, volume label) are written directly
#  - "Price" column (D) cells hold numeric-looking text (e.g. "275.88",
#    "0.00000000750") that must stay literal text (matches the source
#    workbook's inlineStr cells) instead of being auto-coerced to numbers
#    by Excel (which would mangle values like "0.00000000750" into
#    "7.5E-09" and silently reformat trailing/leading zeros). Forcing the
#    cell to a text number-format before assigning the value, then
#    restoring the default "Normal" style, keeps the value an exact text
#    string while leaving the cell's style as it was originally (no
#    explicit style index).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ref, $value) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

function Set-TextValue($ref, $value) {
    $ws.Range($ref).Value = $value
}

# --- Row 2 (BNB): price tick ---
Set-TextCell "D2" "275.88"

# --- Row 3 (OKB): price tick ---
Set-TextCell "D3" "21.16"

# --- Row 4 (HuobiToken): price tick ---
Set-TextCell "D4" "6.267"

# --- Row 5 (Cronos): price tick ---
Set-TextCell "D5" "0.06223"

# --- Row 7: FTXToken -> KuCoinToken (ranking swap with row 8) ---
Set-TextValue "B7" "KuCoinToken"
Set-TextValue "C7" "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextCell  "D7" "6.572"
Set-TextValue "E7" "6KuCoinTokenKCS"

# --- Row 8: KuCoinToken -> FTXToken (ranking swap with row 7) ---
Set-TextValue "B8" "FTXToken"
Set-TextValue "C8" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextCell  "D8" "1.487"
Set-TextValue "E8" "7FTXTokenFTT"

# --- Row 9 (MXToken): price tick ---
Set-TextCell "D9" "0.8280"

# --- Row 10 (WazirX): price tick ---
Set-TextCell "D10" "0.1664"

# --- Row 11 (MandalaExchangeToken): price tick ---
Set-TextCell "D11" "0.08325"

# --- Row 12 (LiechtensteinCryptoassetsExchange): price tick ---
Set-TextCell "D12" "0.03528"

# --- Row 13 (BitrueCoin): price tick ---
Set-TextCell "D13" "0.03167"

# --- Row 14 (BitMartToken): price tick ---
Set-TextCell "D14" "0.09182"

# --- Row 15 (MCDex): price tick ---
Set-TextCell "D15" "3.763"

# --- Row 16 (BitForexToken): price tick ---
Set-TextCell "D16" "0.001630"

# --- Row 17 (CoinExToken): price tick ---
Set-TextCell "D17" "0.04672"

# --- Row 18 (TigerCash): price tick ---
Set-TextCell "D18" "0.006335"

# --- Row 19 (HotbitToken): price tick ---
Set-TextCell "D19" "0.006220"

# --- Row 21 (NitroEx): price tick ---
Set-TextCell "D21" "0.0001500"

# --- Row 22 (LEO): price tick ---
Set-TextCell "D22" "3.720"

# --- Row 23 (BTSEToken): price tick ---
Set-TextCell "D23" "2.322"

# --- Row 24 (One): price tick ---
Set-TextCell "D24" "0.01395"

# --- Row 26 (ProBitToken): price tick ---
Set-TextCell "D26" "0.1242"

# --- Row 40 (IDEX): price tick ---
Set-TextCell "D40" "0.04747"

# --- Row 41: KickToken -> CEJI (ranking shuffle with rows 42/43) ---
Set-TextValue "B41" "CEJI"
Set-TextValue "C41" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextCell  "D41" "0.005198"
Set-TextValue "E41" "40CEJICEJI"

# --- Row 42: BKEXToken -> KickToken (ranking shuffle with rows 41/43) ---
Set-TextValue "B42" "KickToken"
Set-TextValue "C42" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextCell  "D42" "0.007060"
Set-TextValue "E42" "41KickTokenKICKBestin24h"

# --- Row 43: CEJI -> BKEXToken (ranking shuffle with rows 41/42) ---
Set-TextValue "B43" "BKEXToken"
Set-TextValue "C43" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextCell  "D43" "0.1123"
Set-TextValue "E43" "42BKEXTokenBKK"

# --- Row 44 (LocalTraders): price tick ---
Set-TextCell "D44" "0.01145"

# --- Row 45 (CoinLion): price tick ---
Set-TextCell "D45" "0.00006279"

# --- Row 46 (ACDXExchange): price tick + label lost "Bestin24h" flag ---
Set-TextCell  "D46" "0.0009895"
Set-TextValue "E46" "45ACDXExchangeACXT"

# --- Row 47 (Kangarootoken): price tick ---
Set-TextCell "D47" "0.00000000750"

# --- Row 48 (CoinbaseStockToken): price tick ---
Set-TextCell "D48" "0.7347"

# --- Row 49 (BOLO): price tick ---
Set-TextCell "D49" "0.001401"

# --- Row 50 (CryptobidCoin): price tick ---
Set-TextCell "D50" "0.00001899"

# --- Row 51 (SpecialPowerGold): price tick ---
Set-TextCell "D51" "0.01240"
